# JobStatsBurnupChart3: the burn-up chart's sample data window was moved back
# two weeks (14 days) - every date in the "Dates" column (A2:A15) decreases
# by 14, e.g. 5/20/2019 -> 5/6/2019 ... 6/2/2019 -> 5/19/2019. All other
# columns (Story points / Total Points / Completed / Estimated work) and the
# chart itself (still reading Sheet1!$A$2:$A$15 etc.) are left untouched -
# the chart just re-renders against the new dates.
#
# Note: the c:lineChart/c:catAx/c:valAx <c:axId>/<c:crossAx> values are
# internal, engine-minted axis identifiers (Axis.AxisID is read-only in the
# real Excel object model too), so they aren't something user code sets -
# only the underlying cell data is edited here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$startRow = 2
$endRow = 15
$daysBack = 14

for ($row = $startRow; $row -le $endRow; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $cell.Value2 - $daysBack
}
